$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 32
$ws.Range("A3").Value = 52
$ws.Range("A4").Value = 68
$ws.Range("A5").Value = 152
$ws.Range("A7").Value = 188
$ws.Range("A8").Value = 218
$ws.Range("A9").Value = 222
$ws.Range("A10").Value = 320
$ws.Range("A11").Value = 328
$ws.Range("A12").Value = 340
$ws.Range("A13").Value = 558
$ws.Range("A14").Value = 600
$ws.Range("A15").Value = 858
$ws.Range("A16").Value = 862
